# The deck's canonical OOXML ends up with the contents of ppt/theme/theme1.xml
# ("Office Theme") and ppt/theme/theme2.xml ("Integral") swapped: the slide
# master (and therefore the whole deck's applied design) switches from the
# custom "Integral" theme to the stock "Office Theme" palette, while the
# former Integral palette is left behind on the notes master's theme slot.
#
# Re-create that effect through the PowerPoint object model by rewriting the
# active theme's 12 scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) from the Integral values to the stock Office values, in place,
# on ActivePresentation.SlideMaster.Theme.ThemeColorScheme -- the font
# scheme (Arial-based "Office") and format scheme are already identical
# between the two themes, so the color scheme is the only piece that needs
# to change to reproduce the new design.

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme scheme colors, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6,
# hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = BGR($officeColors[$i - 1])
}
